$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.728.46'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").Value = '2.340.12'
$ws.Range("E3").Value = '  +5.40%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '269.98'
$ws.Range("E5").Value = '  -1.50%  '

$ws.Range("D6").Value = '94.52'
$ws.Range("E6").Value = '  +9.04%  '

$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +1.22%  '

$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  +3.11%  '

$ws.Range("D10").Value = '44.89'
$ws.Range("E10").Value = '  +0.24%  '

$ws.Range("D11").Value = '0.0942'
$ws.Range("E11").Value = '  +2.74%  '

$ws.Range("D12").Value = '8.09'
$ws.Range("E12").Value = '  +5.69%  '

$ws.Range("E13").Value = '  +0.43%  '

$ws.Range("D14").Value = '2.679.53'
$ws.Range("E14").Value = '  +4.50%  '

$ws.Range("D15").Value = '15.58'
$ws.Range("E15").Value = '  +5.11%  '

$ws.Range("D16").Value = '0.861'
$ws.Range("E16").Value = '  +9.37%  '

$ws.Range("D17").Value = '2.328.90'

$ws.Range("D18").Value = '43.693.80'
$ws.Range("E18").Value = '  +0.11%  '

$ws.Range("E19").Value = '  +4.85%  '

$ws.Range("D20").Value = '6.38'
$ws.Range("E20").Value = '  +7.49%  '

$ws.Range("D21").Value = '71.72'
$ws.Range("E21").Value = '  +2.43%  '

$ws.Range("D22").Value = '238.16'
$ws.Range("E22").Value = '  +2.81%  '

$ws.Range("D23").Value = '2.26'
$ws.Range("E23").Value = '  -3.30%  '

$ws.Range("D24").Value = '9.63'
$ws.Range("E24").Value = '  +11.43%  '

$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26").Value = '11.35'
$ws.Range("E26").Value = '  +5.76%  '

$ws.Range("D27").Value = '2.51'
$ws.Range("E27").Value = '  -1.26%  '

$ws.Range("E28").Value = '  -1.78%  '

$ws.Range("E29").Value = '  -1.21%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '38.65'
$ws.Range("E30").Value = '  -1.37%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '22.75'
$ws.Range("E31").Value = '  +10.14%  '

$ws.Range("D32").Value = '171.99'
$ws.Range("E32").Value = '  -0.39%  '

$ws.Range("D33").Value = '0.0897'
$ws.Range("E33").Value = '  -0.84%  '

$ws.Range("D34").Value = '5.50'
$ws.Range("E34").Value = '  +3.32%  '

$ws.Range("D35").Value = '0.127'
$ws.Range("E35").Value = '  +3.12%  '

$ws.Range("E36").Value = '  +1.55%  '

$ws.Range("E37").Value = '  -2.47%  '

$ws.Range("D38").Value = '4.37'
$ws.Range("E38").Value = '  +1.78%  '

$ws.Range("D39").Value = '3.41'
$ws.Range("E39").Value = '  +1.37%  '

$ws.Range("E40").Value = '  +9.26%  '

$ws.Range("D41").Value = '0.233'
$ws.Range("E41").Value = '  +13.85%  '

$ws.Range("D42").Value = '1.38'
$ws.Range("E42").Value = '  +22.96%  '

$ws.Range("D43").Value = '12.08'
$ws.Range("E43").Value = '  -2.44%  '

$ws.Range("E44").Value = '  +7.85%  '

$ws.Range("D45").Value = '61.80'
$ws.Range("E45").Value = '  -2.57%  '

$ws.Range("D46").Value = '5.38'
$ws.Range("E46").Value = '  +0.12%  '

$ws.Range("E47").Value = '  +4.12%  '

$ws.Range("D48").Value = '100.64'
$ws.Range("E48").Value = '  +0.92%  '

$ws.Range("D49").Value = '1.21'
$ws.Range("E49").Value = '  +1.39%  '

$ws.Range("D50").Value = '2.562.10'
$ws.Range("E50").Value = '  +4.67%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.364.37'
$ws.Range("E51").Value = '  +4.80%  '
